$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Variavel) text is constant across all data rows
$variavel = "Massa de rendimento mensal real das pessoas de 14 anos ou mais de idade ocupadas na semana de referência com rendimento de trabalho, habitualmente recebido em todos os trabalhos"

# Force column C (Trimestre) to be stored as plain text so values like
# "01/01/2019" are not auto-converted to date serial numbers.
$ws.Range("C2:C82").NumberFormat = "@"

# Data rows: Regiao, Trimestre, Valor (row 1 is the header, already correct)
$data = @(
    @("Brasil", "01/01/2019", 400551),
    @("Brasil", "01/04/2019", 401843),
    @("Brasil", "01/07/2019", 405148),
    @("Brasil", "01/10/2019", 412288),
    @("Brasil", "01/01/2020", 389048),
    @("Brasil", "01/04/2020", 365913),
    @("Brasil", "01/07/2020", 366329),
    @("Brasil", "01/10/2020", 367882),
    @("Brasil", "01/01/2021", 330758),
    @("Brasil", "01/04/2021", 329970),
    @("Brasil", "01/07/2021", 329429),
    @("Brasil", "01/10/2021", 327674),
    @("Brasil", "01/01/2022", 312753),
    @("Brasil", "01/04/2022", 326417),
    @("Brasil", "01/07/2022", 342148),
    @("Brasil", "01/10/2022", 348970),
    @("Brasil", "01/01/2023", 330653),
    @("Brasil", "01/04/2023", 333987),
    @("Brasil", "01/07/2023", 342839),
    @("Brasil", "01/10/2023", 349896),
    @("Brasil", "01/01/2024", 336116),
    @("Brasil", "01/04/2024", 347632),
    @("Brasil", "01/07/2024", 350369),
    @("Brasil", "01/10/2024", 358115),
    @("Brasil", "01/01/2025", 343280),
    @("Brasil", "01/04/2025", 353215),
    @("Brasil", "01/07/2025", 354564),
    @("Nordeste", "01/01/2019", 61993),
    @("Nordeste", "01/04/2019", 62286),
    @("Nordeste", "01/07/2019", 61963),
    @("Nordeste", "01/10/2019", 63533),
    @("Nordeste", "01/01/2020", 59460),
    @("Nordeste", "01/04/2020", 54337),
    @("Nordeste", "01/07/2020", 52489),
    @("Nordeste", "01/10/2020", 54633),
    @("Nordeste", "01/01/2021", 48171),
    @("Nordeste", "01/04/2021", 48650),
    @("Nordeste", "01/07/2021", 49511),
    @("Nordeste", "01/10/2021", 49301),
    @("Nordeste", "01/01/2022", 45989),
    @("Nordeste", "01/04/2022", 48271),
    @("Nordeste", "01/07/2022", 51198),
    @("Nordeste", "01/10/2022", 51906),
    @("Nordeste", "01/01/2023", 49532),
    @("Nordeste", "01/04/2023", 49752),
    @("Nordeste", "01/07/2023", 51000),
    @("Nordeste", "01/10/2023", 51917),
    @("Nordeste", "01/01/2024", 49703),
    @("Nordeste", "01/04/2024", 53508),
    @("Nordeste", "01/07/2024", 53988),
    @("Nordeste", "01/10/2024", 55584),
    @("Nordeste", "01/01/2025", 53053),
    @("Nordeste", "01/04/2025", 54826),
    @("Nordeste", "01/07/2025", 55902),
    @("Sergipe", "01/01/2019", 2643),
    @("Sergipe", "01/04/2019", 2666),
    @("Sergipe", "01/07/2019", 2637),
    @("Sergipe", "01/10/2019", 2684),
    @("Sergipe", "01/01/2020", 2616),
    @("Sergipe", "01/04/2020", 2458),
    @("Sergipe", "01/07/2020", 2277),
    @("Sergipe", "01/10/2020", 2667),
    @("Sergipe", "01/01/2021", 2182),
    @("Sergipe", "01/04/2021", 2347),
    @("Sergipe", "01/07/2021", 2312),
    @("Sergipe", "01/10/2021", 2416),
    @("Sergipe", "01/01/2022", 2126),
    @("Sergipe", "01/04/2022", 2132),
    @("Sergipe", "01/07/2022", 2206),
    @("Sergipe", "01/10/2022", 2302),
    @("Sergipe", "01/01/2023", 2132),
    @("Sergipe", "01/04/2023", 2220),
    @("Sergipe", "01/07/2023", 2186),
    @("Sergipe", "01/10/2023", 2129),
    @("Sergipe", "01/01/2024", 2221),
    @("Sergipe", "01/04/2024", 2349),
    @("Sergipe", "01/07/2024", 2432),
    @("Sergipe", "01/10/2024", 2599),
    @("Sergipe", "01/01/2025", 2385),
    @("Sergipe", "01/04/2025", 2377),
    @("Sergipe", "01/07/2025", 2613),
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $variavel
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $row++
}

Write-Host "Updated $($row - 2) data rows (rows 2..$($row - 1))"
